$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Normal style: set the document's default Western font to Baskerville
#    (adds <w:rPr><w:rFonts w:ascii="Baskerville" w:hAnsi="Baskerville"/></w:rPr>)
# ---------------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.Font.Name = "Baskerville"

# ---------------------------------------------------------------------------
# 2. Heading1Char style: switch the major/heading Latin font from the theme
#    font to an explicit Baskerville (keeps the East Asian / complex-script
#    theme fonts as-is).
# ---------------------------------------------------------------------------
$h1Char = $d.Styles("Heading1Char")
$h1Char.Font.Name = "Baskerville"

# ---------------------------------------------------------------------------
# 3. GeneratedTitle style: also pin it to Baskerville.
# ---------------------------------------------------------------------------
$generatedTitle = $d.Styles("GeneratedTitle")
$generatedTitle.Font.Name = "Baskerville"

# ---------------------------------------------------------------------------
# 4. New styles for margin notes: "marginOuter" (paragraph) linked with
#    "marginOuterChar" (character), both based on the document defaults and
#    set to a small 10pt Baskerville run.
# ---------------------------------------------------------------------------
$marginOuter = $d.Styles.Add("marginOuter", 1)
$marginOuter.NameLocal = "marginOuter"
$marginOuter.BaseStyle = $d.Styles("Normal")
$marginOuter.NextParagraphStyle = $d.Styles("Normal")

$marginOuterChar = $d.Styles.Add("marginOuterChar", 2)
$marginOuterChar.NameLocal = "marginOuter Char"
$marginOuterChar.BaseStyle = $d.Styles("DefaultParagraphFont")

$marginOuter.LinkStyle = $marginOuterChar
$marginOuterChar.LinkStyle = $marginOuter

$marginOuter.QuickStyle = $true

$marginOuter.Font.Size = 10
$marginOuter.Font.SizeBi = 10

$marginOuterChar.Font.Name = "Baskerville"
$marginOuterChar.Font.Size = 10
$marginOuterChar.Font.SizeBi = 10

Write-Host "margin note styles applied"
